# Fruta / hortaliza, semanal
# Inserts 3 new weekly price records for "Femacal de La Calera - Frutilla"
# right before the existing row 270, shifting the previous rows 270-341
# down to 273-344 (dimension grows from A1:T341 to A1:T344).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows above the current row 270 (formatting/styles of
# row 270 - e.g. the date-formatted column D - get carried into the new
# rows automatically, matching the surrounding rows).
$ws.Range("A270:T272").EntireRow.Insert()

# New row 270
$ws.Cells.Item(270, 1).Value = 3
$ws.Cells.Item(270, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(270, 3).Value = "Coquimbo"
$ws.Cells.Item(270, 4).Value = 44841
$ws.Cells.Item(270, 5).Value = 5
$ws.Cells.Item(270, 6).Value = "Fruta"
$ws.Cells.Item(270, 7).Value = 100101
$ws.Cells.Item(270, 8).Value = "Berries"
$ws.Cells.Item(270, 9).Value = 100112025
$ws.Cells.Item(270, 10).Value = "Frutilla"
$ws.Cells.Item(270, 11).Value = "Sin especificar"
$ws.Cells.Item(270, 12).Value = "Especial"
$ws.Cells.Item(270, 13).Value = 58
$ws.Cells.Item(270, 14).Value = 14000
$ws.Cells.Item(270, 15).Value = 14000
$ws.Cells.Item(270, 16).Value = 14000
$ws.Cells.Item(270, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(270, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(270, 19).Value = 2000
$ws.Cells.Item(270, 20).Value = 7

# New row 271
$ws.Cells.Item(271, 1).Value = 3
$ws.Cells.Item(271, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(271, 3).Value = "Coquimbo"
$ws.Cells.Item(271, 4).Value = 44841
$ws.Cells.Item(271, 5).Value = 5
$ws.Cells.Item(271, 6).Value = "Fruta"
$ws.Cells.Item(271, 7).Value = 100101
$ws.Cells.Item(271, 8).Value = "Berries"
$ws.Cells.Item(271, 9).Value = 100112025
$ws.Cells.Item(271, 10).Value = "Frutilla"
$ws.Cells.Item(271, 11).Value = "Sin especificar"
$ws.Cells.Item(271, 12).Value = "Primera"
$ws.Cells.Item(271, 13).Value = 57
$ws.Cells.Item(271, 14).Value = 12000
$ws.Cells.Item(271, 15).Value = 12000
$ws.Cells.Item(271, 16).Value = 12000
$ws.Cells.Item(271, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(271, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(271, 19).Value = 1714
$ws.Cells.Item(271, 20).Value = 7

# New row 272
$ws.Cells.Item(272, 1).Value = 3
$ws.Cells.Item(272, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(272, 3).Value = "Coquimbo"
$ws.Cells.Item(272, 4).Value = 44841
$ws.Cells.Item(272, 5).Value = 5
$ws.Cells.Item(272, 6).Value = "Fruta"
$ws.Cells.Item(272, 7).Value = 100101
$ws.Cells.Item(272, 8).Value = "Berries"
$ws.Cells.Item(272, 9).Value = 100112025
$ws.Cells.Item(272, 10).Value = "Frutilla"
$ws.Cells.Item(272, 11).Value = "Sin especificar"
$ws.Cells.Item(272, 12).Value = "Segunda"
$ws.Cells.Item(272, 13).Value = 40
$ws.Cells.Item(272, 14).Value = 9000
$ws.Cells.Item(272, 15).Value = 9000
$ws.Cells.Item(272, 16).Value = 9000
$ws.Cells.Item(272, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(272, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(272, 19).Value = 1286
$ws.Cells.Item(272, 20).Value = 7
